# Fruta / hortaliza, semanal
# Insert a new weekly price-record row at row 12 (pushing the existing
# rows 12-56 down to 13-57) and populate it with the new data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("12:12").Insert()

$ws.Cells.Item(12, 1).Value  = 7
$ws.Cells.Item(12, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(12, 3).Value  = "Ñuble"
$ws.Cells.Item(12, 4).Value  = 44575
$ws.Cells.Item(12, 5).Value  = 16
$ws.Cells.Item(12, 6).Value  = 100112021
$ws.Cells.Item(12, 7).Value  = "Ají"
$ws.Cells.Item(12, 8).Value  = "Americana (o)"
$ws.Cells.Item(12, 9).Value  = "Primera"
$ws.Cells.Item(12, 10).Value = 80
$ws.Cells.Item(12, 11).Value = 15000
$ws.Cells.Item(12, 12).Value = 16000
$ws.Cells.Item(12, 13).Value = 15500
$ws.Cells.Item(12, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(12, 15).Value = "Región del Maule"
$ws.Cells.Item(12, 16).Value = 1033
$ws.Cells.Item(12, 17).Value = 15
$ws.Cells.Item(12, 18).Value = "Hortaliza"
